# Schulferien - unnötige Zeilen gelöscht
# The legend/footnote block (rows 20-26, columns A:G) contained explanatory
# text ("Stand: ...", source remarks, footnotes) that is no longer needed.
# Clear the cell contents (values) while keeping the existing cell
# formatting/styles and row heights intact - equivalent to selecting the
# range and pressing Delete / "Clear Contents" in Excel.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20:G26").ClearContents()

# Leave the selection where the user ended up after clearing the block.
$ws.Range("A18:G27").Select()
